# Edit script: insert two new price rows (dated 2023-05-28, serial 45041)
# for "Coliflor" at "Vega Modelo de Temuco", pushing the existing rows
# 569-600 down to 571-602.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 569. This shifts the
# existing rows 569:600 down to 571:602, carrying their values/styles
# with them (matches Excel's native Insert behaviour).
$ws.Rows("569:570").Insert()

# New row 569
$ws.Cells.Item(569, 1).Value = 10
$ws.Cells.Item(569, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(569, 3).Value = "La Araucanía"
$ws.Cells.Item(569, 4).Value = 45041
$ws.Cells.Item(569, 5).Value = 9
$ws.Cells.Item(569, 6).Value = 100112008
$ws.Cells.Item(569, 7).Value = "Coliflor"
$ws.Cells.Item(569, 8).Value = "Sin especificar"
$ws.Cells.Item(569, 9).Value = "Primera"
$ws.Cells.Item(569, 10).Value = 500
$ws.Cells.Item(569, 11).Value = 1300
$ws.Cells.Item(569, 12).Value = 1300
$ws.Cells.Item(569, 13).Value = 1300
$ws.Cells.Item(569, 14).Value = "$/unidad"
$ws.Cells.Item(569, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(569, 16).Value = 1300
$ws.Cells.Item(569, 17).Value = 1
$ws.Cells.Item(569, 18).Value = "Hortaliza"

# New row 570
$ws.Cells.Item(570, 1).Value = 10
$ws.Cells.Item(570, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(570, 3).Value = "La Araucanía"
$ws.Cells.Item(570, 4).Value = 45041
$ws.Cells.Item(570, 5).Value = 9
$ws.Cells.Item(570, 6).Value = 100112008
$ws.Cells.Item(570, 7).Value = "Coliflor"
$ws.Cells.Item(570, 8).Value = "Sin especificar"
$ws.Cells.Item(570, 9).Value = "Primera"
$ws.Cells.Item(570, 10).Value = 1500
$ws.Cells.Item(570, 11).Value = 1300
$ws.Cells.Item(570, 12).Value = 1300
$ws.Cells.Item(570, 13).Value = 1300
$ws.Cells.Item(570, 14).Value = "$/unidad"
$ws.Cells.Item(570, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(570, 16).Value = 1300
$ws.Cells.Item(570, 17).Value = 1
$ws.Cells.Item(570, 18).Value = "Hortaliza"
